$d = $word.ActiveDocument

function New-PkgXml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Hunk 1 -----------------------------------------------------------
# Old: para "CASO DE ESTUDIO: COMPLEJO DEPORTIVO U.M.S.S." + empty (sz44) paragraph
# New: an (now empty) title paragraph carrying ind/jc + sz40/es-419, and an
# empty bold sz24/es-419 paragraph. The following MODALIDAD paragraph is
# left untouched.
$pCaso = $d.Paragraphs.Item(10)
$pEmpty44 = $d.Paragraphs.Item(11)
$rng1 = $d.Range($pCaso.Range.Start, $pEmpty44.Range.End)

$body1 = '<w:body>' +
  '<w:p><w:pPr><w:ind w:left="709" w:hanging="1"/><w:jc w:val="center"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p>' +
  '</w:body>'
$rng1.InsertXML((New-PkgXml $body1))

# --- Hunk 2 -------------------------------------------------------------
# Old: TUTOR paragraph (ends with bookmarkStart/bookmarkEnd) + empty
# paragraph + "PERIODO II - 2015" paragraph (3 runs) + trailing empty <w:p/>.
# New: TUTOR paragraph (bookmark removed) + empty paragraph + two new empty
# center-indented paragraphs + a paragraph holding the bookmark and a single
# merged "PERIODO II - 2015" run. The trailing empty paragraph disappears.
$pTutor = $d.Paragraphs.Item(16)
$pTrailingEmpty = $d.Paragraphs.Item(19)
$rng2 = $d.Range($pTutor.Range.Start, $pTrailingEmpty.Range.End)

$body2 = '<w:body>' +
  '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="005E228D"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t>TUTOR:</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> Lic. </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t>Valentin</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t>Laime</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> Zapata</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:ind w:left="708" w:firstLine="708"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:ind w:left="708" w:firstLine="708"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:ind w:left="708" w:firstLine="708"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-419"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr><w:t>PERIODO II - 2015</w:t></w:r>' +
  '</w:p>' +
  '</w:body>'
$rng2.InsertXML((New-PkgXml $body2))

Write-Output "done"
